# Applies the two changes described by the commit:
#  1. Slide 5's table switches from the custom "Table_0" style to the
#     built-in table style {E4674AB0-7CF2-4FCF-B325-6C42FD35AA48}.
#  2. The deck's design theme (ppt/theme/theme1.xml, the part bound to the
#     slide master) is recoloured from the "Integral" / "Red Violet" scheme
#     to the stock Office theme palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on Slide 5 --------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{E4674AB0-7CF2-4FCF-B325-6C42FD35AA48}")

# --- 2. Swap the presentation theme colours to the Office theme -----------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.ColorScheme

$colorScheme.Colors(1).RGB  = 0          # dk1      -> 000000
$colorScheme.Colors(2).RGB  = 16777215   # lt1      -> FFFFFF
$colorScheme.Colors(3).RGB  = 6968388    # dk2      -> 44546A
$colorScheme.Colors(4).RGB  = 15132391   # lt2      -> E7E6E6
$colorScheme.Colors(5).RGB  = 13998939   # accent1  -> 5B9BD5
$colorScheme.Colors(6).RGB  = 3243501    # accent2  -> ED7D31
$colorScheme.Colors(7).RGB  = 10855845   # accent3  -> A5A5A5
$colorScheme.Colors(8).RGB  = 49407      # accent4  -> FFC000
$colorScheme.Colors(9).RGB  = 12874308   # accent5  -> 4472C4
$colorScheme.Colors(10).RGB = 4697456    # accent6  -> 70AD47
$colorScheme.Colors(11).RGB = 12673797   # hlink    -> 0563C1
$colorScheme.Colors(12).RGB = 7491477    # folHlink -> 954F72
